{"js": "// Replace each \"NN\u00f7N=\" equation in the document's table with a new one,\n// in document order. Some original equations repeat (e.g. \"12\u00f78=\" and\n// \"40\u00f73=\" each appear twice) and each occurrence maps to a different\n// replacement, so we cannot do a blind \"replace all\" per search term \u2014\n// instead we walk the ordered list of (oldText, newText) pairs and, for\n// each one, replace only the first still-matching occurrence left in the\n// document. Because replacements never reintroduce an old value that a\n// later pair still needs to match, processing in this fixed order is safe.\nconst replacements = [\n  [\"13\u00f78=\", \"76\u00f76=\"],\n  [\"53\u00f72=\", \"15\u00f78=\"],\n  [\"69\u00f73=\", \"20\u00f78=\"],\n  [\"51\u00f78=\", \"82\u00f74=\"],\n  [\"22\u00f73=\", \"52\u00f74=\"],\n  [\"18\u00f72=\", \"61\u00f72=\"],\n  [\"16\u00f74=\", \"29\u00f79=\"],\n  [\"12\u00f78=\", \"94\u00f77=\"],\n  [\"68\u00f76=\", \"45\u00f76=\"],\n  [\"92\u00f79=\", \"43\u00f73=\"],\n  [\"44\u00f72=\", \"88\u00f75=\"],\n  [\"40\u00f73=\", \"47\u00f72=\"],\n  [\"61\u00f78=\", \"87\u00f73=\"],\n  [\"76\u00f79=\", \"13\u00f79=\"],\n  [\"57\u00f73=\", \"41\u00f79=\"],\n  [\"18\u00f73=\", \"25\u00f73=\"],\n  [\"63\u00f72=\", \"39\u00f78=\"],\n  [\"24\u00f75=\", \"46\u00f78=\"],\n  [\"65\u00f75=\", \"27\u00f78=\"],\n  [\"12\u00f78=\", \"56\u00f75=\"],\n  [\"54\u00f75=\", \"51\u00f75=\"],\n  [\"55\u00f72=\", \"72\u00f77=\"],\n  [\"40\u00f73=\", \"76\u00f73=\"],\n  [\"36\u00f79=\", \"17\u00f79=\"],\n  [\"64\u00f75=\", \"55\u00f77=\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: \"${oldText}\"`);\n  }\n\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Replace each \"NN\u00f7N=\" equation in the document's table with a new one,\n# in document order. Some original equations repeat (e.g. \"12\u00f78=\" and\n# \"40\u00f73=\" each appear twice) and each occurrence maps to a different\n# replacement, so a single MatchCase Find/Replace-All per term would not\n# be correct. Instead we walk the ordered list of (oldText, newText)\n# pairs and, for each one, re-run Find.Execute with Replace = wdReplaceOne\n# (1) starting from the top of the document each time; this always lands\n# on the first occurrence that still has its original text, since earlier\n# pairs have already been rewritten.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"13\u00f78=\", \"76\u00f76=\"),\n  @(\"53\u00f72=\", \"15\u00f78=\"),\n  @(\"69\u00f73=\", \"20\u00f78=\"),\n  @(\"51\u00f78=\", \"82\u00f74=\"),\n  @(\"22\u00f73=\", \"52\u00f74=\"),\n  @(\"18\u00f72=\", \"61\u00f72=\"),\n  @(\"16\u00f74=\", \"29\u00f79=\"),\n  @(\"12\u00f78=\", \"94\u00f77=\"),\n  @(\"68\u00f76=\", \"45\u00f76=\"),\n  @(\"92\u00f79=\", \"43\u00f73=\"),\n  @(\"44\u00f72=\", \"88\u00f75=\"),\n  @(\"40\u00f73=\", \"47\u00f72=\"),\n  @(\"61\u00f78=\", \"87\u00f73=\"),\n  @(\"76\u00f79=\", \"13\u00f79=\"),\n  @(\"57\u00f73=\", \"41\u00f79=\"),\n  @(\"18\u00f73=\", \"25\u00f73=\"),\n  @(\"63\u00f72=\", \"39\u00f78=\"),\n  @(\"24\u00f75=\", \"46\u00f78=\"),\n  @(\"65\u00f75=\", \"27\u00f78=\"),\n  @(\"12\u00f78=\", \"56\u00f75=\"),\n  @(\"54\u00f75=\", \"51\u00f75=\"),\n  @(\"55\u00f72=\", \"72\u00f77=\"),\n  @(\"40\u00f73=\", \"76\u00f73=\"),\n  @(\"36\u00f79=\", \"17\u00f79=\"),\n  @(\"64\u00f75=\", \"55\u00f77=\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $range = $d.Content\n  $find = $range.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n\n  # FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n  # MatchAllWordForms, Forward, Wrap(wdFindContinue=1), Format,\n  # ReplaceWith, Replace(wdReplaceOne=1)\n  $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 1)\n\n  if (-not $found) {\n    throw \"Could not find text to replace: $oldText\"\n  }\n}\n"}
